$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Refresh the "as-of" timestamp (column D) for all existing data rows (2-58) ----
$ws.Range("D2:D58").Value = 45960.292199074072

# ---- 2) Replace the refreshed portion of the report (rows 19-58) with the new pull ----
$data1958 = @(
    ,@(19, "长沙特来电飞狐四方坪西区充电站", "604号直流", 45957.218495370369)
    ,@(20, "长沙特来电飞狐四方坪南区充电站", "406号直流", 45957.294004629628)
    ,@(21, "长沙特来电飞狐四方坪西区充电站", "702号直流", 45958.053842592592)
    ,@(22, "长沙特来电飞狐四方坪西区充电站", "602号直流", 45958.233749999999)
    ,@(23, "长沙特来电飞狐四方坪西区充电站", "903号直流", 45958.509386574071)
    ,@(24, "长沙特来电飞狐四方坪南区充电站", "105号直流", 45958.544050925928)
    ,@(25, "长沙特来电飞狐四方坪东区充电站", "004A号直流", 45958.647453703707)
    ,@(26, "长沙特来电飞狐四方坪西区充电站", "603号直流", 45959.031655092593)
    ,@(27, "长沙特来电飞狐四方坪西区充电站", "305号直流", 45959.041307870371)
    ,@(28, "长沙特来电飞狐四方坪西区充电站", "904号直流", 45959.042893518519)
    ,@(29, "长沙特来电飞狐四方坪东区充电站", "905号直流", 45959.052430555559)
    ,@(30, "长沙特来电飞狐四方坪西区充电站", "801号直流", 45959.071608796294)
    ,@(31, "长沙特来电飞狐四方坪西区充电站", "705号直流", 45959.245324074072)
    ,@(32, "长沙特来电飞狐四方坪南区充电站", "204号直流", 45959.321967592594)
    ,@(33, "长沙市开福区高岭香江国际城充电站建设项目", "105号直流", 45959.364618055559)
    ,@(34, "长沙特来电飞狐四方坪西区充电站", "A03号直流", 45959.368437500001)
    ,@(35, "长沙市开福区高岭香江国际城充电站建设项目", "107号直流", 45959.381921296299)
    ,@(36, "长沙特来电飞狐四方坪西区充电站", "B01号直流", 45959.49728009259)
    ,@(37, "长沙市开福区高岭香江国际城充电站建设项目", "210号直流", 45959.509259259263)
    ,@(38, "长沙特来电飞狐四方坪西区充电站", "405号直流", 45959.525266203702)
    ,@(39, "长沙特来电飞狐四方坪东区充电站", "103号直流", 45959.530555555553)
    ,@(40, "长沙特来电飞狐四方坪南区充电站", "403号直流", 45959.541435185187)
    ,@(41, "长沙特来电飞狐四方坪西区充电站", "A01号直流", 45959.548055555555)
    ,@(42, "长沙特来电飞狐四方坪南区充电站", "104号直流", 45959.550138888888)
    ,@(43, "长沙特来电飞狐四方坪东区充电站", "008B号直流", 45959.55945601852)
    ,@(44, "长沙特来电飞狐四方坪东区充电站", "102号直流", 45959.563090277778)
    ,@(45, "长沙特来电飞狐四方坪南区充电站", "201号直流", 45959.565925925926)
    ,@(46, "长沙特来电飞狐四方坪西区充电站", "B02号直流", 45959.567187499997)
    ,@(47, "长沙特来电飞狐四方坪西区充电站", "401号直流", 45959.571030092593)
    ,@(48, "长沙特来电飞狐四方坪东区充电站", "402号直流", 45959.571319444447)
    ,@(49, "长沙特来电飞狐四方坪西区充电站", "901号直流", 45959.572048611109)
    ,@(50, "长沙特来电飞狐四方坪南区充电站", "203号直流", 45959.58021990741)
    ,@(51, "长沙特来电飞狐四方坪南区充电站", "103号直流", 45959.580601851849)
    ,@(52, "长沙市开福区高岭香江国际城充电站建设项目", "110号直流", 45959.581030092595)
    ,@(53, "长沙市开福区高岭香江国际城充电站建设项目", "101号直流", 45959.585370370369)
    ,@(54, "长沙市开福区高岭香江国际城充电站建设项目", "311号直流", 45959.589849537035)
    ,@(55, "长沙市开福区高岭香江国际城充电站建设项目", "306号直流", 45959.599074074074)
    ,@(56, "长沙特来电飞狐四方坪东区充电站", "001B号直流", 45959.601956018516)
    ,@(57, "长沙特来电飞狐四方坪西区充电站", "A02号直流", 45959.603078703702)
    ,@(58, "长沙市开福区高岭香江国际城充电站建设项目", "108号直流", 45959.640196759261)
)
foreach ($row in $data1958) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# ---- 3) Append newly-observed rows (59-65), cloning formatting from the last existing row ----
$ws.Range("A58:E58").Copy()
$ws.Range("A59:E65").PasteSpecial(-4122)

$data5965 = @(
    ,@(59, "长沙特来电飞狐四方坪南区充电站", "401号直流", 45959.640868055554)
    ,@(60, "长沙特来电飞狐四方坪南区充电站", "301号直流", 45959.643657407411)
    ,@(61, "长沙市开福区高岭香江国际城充电站建设项目", "106号直流", 45959.644641203704)
    ,@(62, "长沙特来电飞狐四方坪西区充电站", "704号直流", 45959.69023148148)
    ,@(63, "长沙特来电飞狐四方坪西区充电站", "804号直流", 45959.692719907405)
    ,@(64, "长沙市开福区高岭香江国际城充电站建设项目", "109号直流", 45959.707060185188)
    ,@(65, "长沙市开福区高岭香江国际城充电站建设项目", "102号直流", 45959.768240740741)
)
foreach ($row in $data5965) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = 45960.292199074072
}

# ---- 4) Update the view state: scroll position + active selection ----
$ws.Range("J49").Select()
try { $excel.ActiveWindow.ScrollRow = 34 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}

